$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 185.22223
$ws.Range("J9").Value = 207.14285
$ws.Range("L9").Value = 207.14285
$ws.Range("N9").Value = -545.14285

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 326.76923
$ws.Range("I38").Value = 326.76923
$ws.Range("K38").Value = 980.30769
$ws.Range("M38").Value = -608.30769

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4922.231
$ws.Range("I40").Value = 4250
$ws.Range("J40").Value = 5044.4546
$ws.Range("K40").Value = 4250
$ws.Range("L40").Value = 5044.4546
$ws.Range("M40").Value = -4075
$ws.Range("N40").Value = -5394.4546

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2146.8572
$ws.Range("I86").Value = 2084.111
$ws.Range("J86").Value = 2259.8
$ws.Range("K86").Value = 2084.111
$ws.Range("L86").Value = 2259.8
$ws.Range("M86").Value = -961.1109999999999
$ws.Range("N86").Value = -4505.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2146.8572
$ws.Range("I89").Value = 2084.111
$ws.Range("J89").Value = 2259.8
$ws.Range("K89").Value = 10420.555
$ws.Range("L89").Value = 11299
$ws.Range("M89").Value = -4804.555
$ws.Range("N89").Value = -22531

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4949.0557
$ws.Range("I137").Value = 2056
$ws.Range("J137").Value = 10067.538
$ws.Range("K137").Value = 6168
$ws.Range("L137").Value = 30202.614
$ws.Range("M137").Value = -3618
$ws.Range("N137").Value = -35302.614

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1521.6818
$ws.Range("I138").Value = 1327.4762
$ws.Range("J138").Value = 5600
$ws.Range("K138").Value = 3982.4286
$ws.Range("L138").Value = 16800
$ws.Range("M138").Value = 1157.5714
$ws.Range("N138").Value = -27080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 11863.333
$ws.Range("J37").Value = 30000
$ws.Range("L37").Value = 30000
$ws.Range("N37").Value = -30546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 15096.392
$ws.Range("I45").Value = 12683.111
$ws.Range("J45").Value = 23784.2
$ws.Range("K45").Value = 12683.111
$ws.Range("L45").Value = 23784.2
$ws.Range("M45").Value = -12306.111
$ws.Range("N45").Value = -24538.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3634.8462
$ws.Range("I61").Value = 3434.913
$ws.Range("K61").Value = 3434.913
$ws.Range("M61").Value = -3222.913

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3069.9412
$ws.Range("I74").Value = 1730
$ws.Range("K74").Value = 1730
$ws.Range("M74").Value = -856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3069.9412
$ws.Range("I77").Value = 1730
$ws.Range("K77").Value = 8650
$ws.Range("M77").Value = -4282

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3634.8462
$ws.Range("I136").Value = 3434.913
$ws.Range("K136").Value = 10304.739
$ws.Range("M136").Value = -7754.739

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 30000
$ws.Range("J35").Value = 30000
$ws.Range("L35").Value = 30000
$ws.Range("N35").Value = -30620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 14839.473
$ws.Range("I105").Value = 3575.2693
$ws.Range("J105").Value = 44126.4
$ws.Range("K105").Value = 3575.2693
$ws.Range("L105").Value = 44126.4
$ws.Range("M105").Value = -1828.2693
$ws.Range("N105").Value = -47620.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 200000220
$ws.Range("I7").Value = 250000210
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 250000210
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -250000097
$ws.Range("N7").Value = -526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1018.2
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 674.75
$ws.Range("I105").Value = 633
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 633
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = 1114
$ws.Range("N105").Value = -4294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1132.1538
$ws.Range("I134").Value = 1022.25
$ws.Range("J134").Value = 1308
$ws.Range("K134").Value = 3066.75
$ws.Range("L134").Value = 3924
$ws.Range("M134").Value = -531.75
$ws.Range("N134").Value = -8994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1018.2
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 911.4
$ws.Range("I5").Value = 707.6316
$ws.Range("J5").Value = 1153.375
$ws.Range("K5").Value = 2122.8948
$ws.Range("L5").Value = 3460.125
$ws.Range("M5").Value = -2010.8948
$ws.Range("N5").Value = -3684.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 111360.89
$ws.Range("I34").Value = 125156
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 375468
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -375384
$ws.Range("N34").Value = -3168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 72799.86
$ws.Range("I39").Value = 100769.9
$ws.Range("J39").Value = 2874.75
$ws.Range("K39").Value = 302309.7
$ws.Range("L39").Value = 8624.25
$ws.Range("M39").Value = -302015.7
$ws.Range("N39").Value = -9212.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 8402560
$ws.Range("J55").Value = 10419750
$ws.Range("L55").Value = 31259250
$ws.Range("N55").Value = -31259604

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2339
$ws.Range("I64").Value = 122
$ws.Range("J64").Value = 3078
$ws.Range("K64").Value = 366
$ws.Range("L64").Value = 9234
$ws.Range("M64").Value = -96
$ws.Range("N64").Value = -9774

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2339
$ws.Range("I67").Value = 122
$ws.Range("J67").Value = 3078
$ws.Range("K67").Value = 366
$ws.Range("L67").Value = 9234
$ws.Range("M67").Value = 570
$ws.Range("N67").Value = -11106

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 710.8
$ws.Range("I107").Value = 790.1667
$ws.Range("J107").Value = 685.7368
$ws.Range("K107").Value = 2370.5001
$ws.Range("L107").Value = 2057.2104
$ws.Range("M107").Value = -450.5001000000002
$ws.Range("N107").Value = -5897.2104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 4460.8887
$ws.Range("I118").Value = 191.16667
$ws.Range("J118").Value = 13000.333
$ws.Range("K118").Value = 573.50001
$ws.Range("L118").Value = 39000.999
$ws.Range("M118").Value = 669.49999
$ws.Range("N118").Value = -41486.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 295.3871
$ws.Range("I122").Value = 214.92857
$ws.Range("J122").Value = 361.64706
$ws.Range("K122").Value = 1934.35713
$ws.Range("L122").Value = 3254.82354
$ws.Range("M122").Value = 515.6428699999999
$ws.Range("N122").Value = -8154.82354

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 911.4
$ws.Range("I135").Value = 707.6316
$ws.Range("J135").Value = 1153.375
$ws.Range("K135").Value = 6368.6844
$ws.Range("L135").Value = 10380.375
$ws.Range("M135").Value = -3833.6844
$ws.Range("N135").Value = -15450.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 36495
$ws.Range("J52").Value = 36495
$ws.Range("L52").Value = 36495
$ws.Range("N52").Value = -37013

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8465.611000000001
$ws.Range("I80").Value = 12099.4
$ws.Range("K80").Value = 12099.4
$ws.Range("M80").Value = -11101.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 8465.611000000001
$ws.Range("I83").Value = 12099.4
$ws.Range("K83").Value = 60497
$ws.Range("M83").Value = -55505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 7442.091
$ws.Range("I99").Value = 2371.6667
$ws.Range("J99").Value = 30259
$ws.Range("K99").Value = 2371.6667
$ws.Range("L99").Value = 30259
$ws.Range("M99").Value = -125.6667000000002
$ws.Range("N99").Value = -34751

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9704.305
$ws.Range("I7").Value = 7849.375
$ws.Range("J7").Value = 13944.143
$ws.Range("K7").Value = 7849.375
$ws.Range("L7").Value = 13944.143
$ws.Range("M7").Value = -7737.375
$ws.Range("N7").Value = -14168.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1409.25
$ws.Range("J16").Value = 850
$ws.Range("L16").Value = 850
$ws.Range("N16").Value = -1190

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1237.3182
$ws.Range("I22").Value = 756.6667
$ws.Range("J22").Value = 1570.0769
$ws.Range("K22").Value = 756.6667
$ws.Range("L22").Value = 1570.0769
$ws.Range("M22").Value = -461.6667
$ws.Range("N22").Value = -2160.0769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1237.3182
$ws.Range("I27").Value = 756.6667
$ws.Range("J27").Value = 1570.0769
$ws.Range("K27").Value = 756.6667
$ws.Range("L27").Value = 1570.0769
$ws.Range("M27").Value = -649.6667
$ws.Range("N27").Value = -1784.0769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4509
$ws.Range("I82").Value = 4737.625
$ws.Range("K82").Value = 4737.625
$ws.Range("M82").Value = -4376.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 4509
$ws.Range("I85").Value = 4737.625
$ws.Range("K85").Value = 4737.625
$ws.Range("M85").Value = -3489.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 25189
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H91").Value = 25189
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4220.769
$ws.Range("I122").Value = 3488
$ws.Range("J122").Value = 4848.857
$ws.Range("K122").Value = 10464
$ws.Range("L122").Value = 14546.571
$ws.Range("M122").Value = -8014
$ws.Range("N122").Value = -19446.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 9704.305
$ws.Range("I126").Value = 7849.375
$ws.Range("J126").Value = 13944.143
$ws.Range("K126").Value = 23548.125
$ws.Range("L126").Value = 41832.429
$ws.Range("M126").Value = -21078.125
$ws.Range("N126").Value = -46772.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5212.3184
$ws.Range("I132").Value = 5111.7334
$ws.Range("J132").Value = 5427.857
$ws.Range("K132").Value = 15335.2002
$ws.Range("L132").Value = 16283.571
$ws.Range("M132").Value = -12805.2002
$ws.Range("N132").Value = -21343.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1548.5
$ws.Range("I122").Value = 1134.4849
$ws.Range("K122").Value = 3403.4547
$ws.Range("M122").Value = -953.4546999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1279701.8
$ws.Range("I132").Value = 1585836.5
$ws.Range("K132").Value = 4757509.5
$ws.Range("M132").Value = -4754979.5
